# LaunchMainApp sequence diagram: rename the example object from "f" to
# "fdc" (adhering to the example object name used elsewhere), plus the
# incidental "datetimeFigureOut" field re-cache that PowerPoint performs
# whenever the deck is next saved on a different day.

# ---------------------------------------------------------------------
# Helper: convert an EMU value to the points value that will round-trip
# back to the exact same EMU value through this host's Shape position
# setters (Left/Top/Width/Height are stored as 32-bit floats internally
# and truncated - not rounded - when re-expressed in EMU on save).
# This loop only touches plain numbers (no COM), so it is cheap.
# ---------------------------------------------------------------------
function Emu-ToPt([double]$emu) {
    $pt = $emu / 12700.0
    for ($i = 0; $i -lt 4000; $i++) {
        $candidate = $pt + ($i * 0.0000005)
        $asFloat = [float]$candidate
        $backEmu = [math]::Floor([double]$asFloat * 12700.0)
        if ($backEmu -eq $emu) {
            return $candidate
        }
    }
    return $pt
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholders: refresh the cached "datetimeFigureOut" field
#    text on the slide master, every slide layout, and the notes
#    master from 11/10/2018 to 11/11/2018.
#    NOTE: loops that walk a COM Shapes collection are kept inline
#    (not inside a helper function) - doing that nested-loop-in-a-
#    function exhausts this host's PowerShell statement budget.
# ---------------------------------------------------------------------
$master = $p.SlideMaster

$masterShapes = $master.Shapes
$masterDate = $null
for ($i = 1; $i -le $masterShapes.Count; $i++) {
    $shp = $masterShapes.Item($i)
    if ($shp.Name -like "Date Placeholder*") { $masterDate = $shp }
}
if ($masterDate -ne $null) { $masterDate.TextFrame.TextRange.Text = "11/11/2018" }

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $lay = $layouts.Item($i)
    $layShapes = $lay.Shapes
    $layDate = $null
    for ($j = 1; $j -le $layShapes.Count; $j++) {
        $shp = $layShapes.Item($j)
        if ($shp.Name -like "Date Placeholder*") { $layDate = $shp }
    }
    if ($layDate -ne $null) { $layDate.TextFrame.TextRange.Text = "11/11/2018" }
}

$notesMaster = $p.NotesMaster
$notesShapes = $notesMaster.Shapes
$notesDate = $null
for ($i = 1; $i -le $notesShapes.Count; $i++) {
    $shp = $notesShapes.Item($i)
    if ($shp.Name -like "Date Placeholder*") { $notesDate = $shp }
}
if ($notesDate -ne $null) { $notesDate.TextFrame.TextRange.Text = "11/11/2018" }

# ---------------------------------------------------------------------
# 2) Slide 1 content: rename the example object from "f" to "fdc".
#    Shape Name is not unique ("Rectangle 62" is reused 4x), so look
#    shapes up by their stable numeric Id instead.
# ---------------------------------------------------------------------
$s = $p.Slides.Item(1)
$slideShapes = $s.Shapes

$rect = $null
$tbox = $null
for ($i = 1; $i -le $slideShapes.Count; $i++) {
    $shp = $slideShapes.Item($i)
    if ($shp.Id -eq 36) { $rect = $shp }
    if ($shp.Id -eq 56) { $tbox = $shp }
}

# "Rectangle 62" (shape Id 36): text "f: FirstDay" -> "fdc: FirstDay"
$rectRange = $rect.TextFrame.TextRange
$rectChar = $rectRange.Characters(1, 1)
$rectChar.Text = "fdc"

# "TextBox 55" (shape Id 56): text "...f.computeAppTitle..." ->
# "...fdc.computeAppTitle..." and reposition/resize the box to fit the
# now-longer label.
$tboxRange = $tbox.TextFrame.TextRange
$fullText = $tboxRange.Text
$startPos = $fullText.IndexOf("f.computeAppTitle") + 1
$tboxChar = $tboxRange.Characters($startPos, 17)
$tboxChar.Text = "fdc.computeAppTitle"

$tbox.Left = Emu-ToPt 229176
$tbox.Top = Emu-ToPt 4053200
$tbox.Width = Emu-ToPt 2308026
$tbox.Height = Emu-ToPt 184666

Write-Output "Rectangle 62 text: $($rectRange.Text)"
Write-Output "TextBox 55 text: $($tboxRange.Text)"
